$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.57%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.55%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.097"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.70%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07944"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.36%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.250"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-7.71%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.771"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.11%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.860"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.22%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.01%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1733"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.05%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.20%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09262"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.24%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03037"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.17%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001504"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006039"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.65%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.476"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.67%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.268"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.44%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.33%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1308"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.43%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.905"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-16.22%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'8.46%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04607"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.10%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.51%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004472"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.87%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003391"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'23.58%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01759"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.11%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04603"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006980"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-5.54%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1360"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.29%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002186"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.27%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009569"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-12.21%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006315"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.69%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000748"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.21%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.007964"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-19.44%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7471"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-8.95%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.21%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("E50").Style = "Normal"
